# Updates cryptocurrency price/volume figures in the "cryptos" sheet.
# D-column "Price" values must stay plain text (they mix thousands-dot
# formatting like "29.390.92" with decimal-looking text like "1.001"), so
# a direct $cell.Value = "1.001" would be auto-coerced by Excel into the
# NUMBER 1.001 (dropping trailing zeros / using scientific notation for
# tiny values, and silently switching the stored type away from text).
# To force literal text without touching any cell style/number-format,
# we stage the text via a `="..."` formula, then Copy + PasteSpecial
# (values only) back onto itself - this collapses the formula down to a
# plain inline/shared string identical in shape to the original cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $escaped = $text.Replace('"', '""')
    $range = $ws.Range($cellRef)
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# Row 2
Set-TextValue "D2" "29.390.92"
$ws.Range("E2").Value = "  +0.28%  "

# Row 3
Set-TextValue "D3" "1.878.23"
$ws.Range("E3").Value = "  +0.27%  "

# Row 4
Set-TextValue "D4" "1.001"
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
Set-TextValue "D5" "0.7166"
$ws.Range("E5").Value = "  +1.19%  "

# Row 6
Set-TextValue "D6" "243.45"
$ws.Range("E6").Value = "  +0.66%  "

# Row 7
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
Set-TextValue "D8" "0.07925"
$ws.Range("E8").Value = "  +1.70%  "

# Row 9
$ws.Range("E9").Value = "  +1.11%  "

# Row 10
Set-TextValue "D10" "24.97"
$ws.Range("E10").Value = "  -0.41%  "

# Row 11
Set-TextValue "D11" "0.08133"
$ws.Range("E11").Value = "  -2.91%  "

# Row 12
Set-TextValue "D12" "1.889.63"
$ws.Range("E12").Value = "  +0.61%  "

# Row 13
Set-TextValue "D13" "94.98"
$ws.Range("E13").Value = "  +4.09%  "

# Row 14
Set-TextValue "D14" "5.234"
$ws.Range("E14").Value = "  -0.01%  "

# Row 15
Set-TextValue "D15" "0.7073"
$ws.Range("E15").Value = "  -1.31%  "

# Row 17
Set-TextValue "D17" "0.000008403"
$ws.Range("E17").Value = "  +0.22%  "

# Row 18
Set-TextValue "D18" "29.409.22"
$ws.Range("E18").Value = "  +0.33%  "

# Row 19
Set-TextValue "D19" "252.10"
$ws.Range("E19").Value = "  +4.95%  "

# Row 20
Set-TextValue "D20" "13.34"
$ws.Range("E20").Value = "  +1.09%  "

# Row 21
Set-TextValue "D21" "2.143.11"
$ws.Range("E21").Value = "  +1.23%  "

# Row 23
Set-TextValue "D23" "7.668"
$ws.Range("E23").Value = "  -0.88%  "

# Row 24
Set-TextValue "D24" "1.001"
$ws.Range("E24").Value = "  +0.02%  "

# Row 25
Set-TextValue "D25" "0.1584"
$ws.Range("E25").Value = "  -0.82%  "

# Row 26
Set-TextValue "D26" "9.066"
$ws.Range("E26").Value = "  +0.35%  "

# Row 27
Set-TextValue "D27" "162.17"
$ws.Range("E27").Value = "  -0.30%  "

# Row 28
Set-TextValue "D28" "18.92"
$ws.Range("E28").Value = "  +2.43%  "

# Row 29
Set-TextValue "D29" "1.507"
$ws.Range("E29").Value = "  +0.17%  "

# Row 30
Set-TextValue "D30" "4.410"
$ws.Range("E30").Value = "  +0.10%  "

# Row 31
Set-TextValue "D31" "4.294"
$ws.Range("E31").Value = "  -1.03%  "

# Row 32
Set-TextValue "D32" "1.219"
$ws.Range("E32").Value = "  +0.48%  "

# Row 33
Set-TextValue "D33" "0.05331"
$ws.Range("E33").Value = "  -0.45%  "

# Row 34
Set-TextValue "D34" "1.945"
$ws.Range("E34").Value = "  +0.25%  "

# Row 35
Set-TextValue "D35" "0.7581"
$ws.Range("E35").Value = "  +1.67%  "

# Row 36
$ws.Range("E36").Value = "  +0.24%  "

# Row 37
Set-TextValue "D37" "2.701"
$ws.Range("E37").Value = "  +0.64%  "

# Row 39
Set-TextValue "D39" "1.270.21"
$ws.Range("E39").Value = "  +2.46%  "

# Row 40
Set-TextValue "D40" "2.760"
$ws.Range("E40").Value = "  +1.04%  "

# Row 41
Set-TextValue "D41" "6.397"
$ws.Range("E41").Value = "  -1.72%  "

# Row 42
Set-TextValue "D42" "112.23"
$ws.Range("E42").Value = "  +2.30%  "

# Row 43
Set-TextValue "D43" "0.9050"
$ws.Range("E43").Value = "  +1.39%  "

# Row 44
Set-TextValue "D44" "74.09"
$ws.Range("E44").Value = "  +2.59%  "

# Row 45
Set-TextValue "D45" "1.001"
$ws.Range("E45").Value = "  +0.10%  "

# Row 46
$ws.Range("E46").Value = "  -0.33%  "

# Row 47
Set-TextValue "D47" "2.039.01"
$ws.Range("E47").Value = "  +1.00%  "

# Row 48
Set-TextValue "D48" "1.809"
$ws.Range("E48").Value = "  +0.93%  "

# Row 49
Set-TextValue "D49" "0.5204"
$ws.Range("E49").Value = "  +0.19%  "

# Row 50
Set-TextValue "D50" "9.512"
$ws.Range("E50").Value = "  +0.69%  "

# Row 51
Set-TextValue "D51" "0.4342"
$ws.Range("E51").Value = "  +0.23%  "

